$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the first data row (row 2): NAIARA -> MARCUS with new account/value
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "005622794"
$ws.Cells.Item(2, 2).Value = "MARCUS"
$ws.Cells.Item(2, 3).Value = 690000

# 2. Insert a new row above row 5 (CAIO) for the new THEOMAR entry
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004231509"
$ws.Cells.Item(5, 2).Value = "THEOMAR"
$ws.Cells.Item(5, 3).Value = 20345.86

# 3. Remove the old MARCUS row (004575621 / 17000), now at row 7
$ws.Rows.Item(7).Delete()

# 4. Remove the VENIA row (004813166 / 7511.38), now at row 8
$ws.Rows.Item(8).Delete()
